$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Active")

# Insert a new row at position 10 (shifts existing rows 10-21 down to 11-22)
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row 10 with the new todo item
$ws.Cells.Item(10, 1).Value = 34
$ws.Cells.Item(10, 2).Value = "when zooming, if a scroll bar is all the way to min or max, keep it there"
$ws.Cells.Item(10, 3).Value = "Todo"
$ws.Cells.Item(10, 4).Value = "Task"

# Created date is stored as literal text (not an auto-converted date serial),
# matching the rest of the "Created" column - force text format first.
$created = $ws.Cells.Item(10, 5)
$created.NumberFormat = "@"
$created.Value = "8/11/2018"
$created.Style = "Normal"

# Bump the Max Id tracker on the Config sheet
$cfg = $wb.Worksheets.Item("Config")
$cfg.Cells.Item(2, 6).Value = 34
